$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13. This shifts the former rows 13-41 down to 14-42,
# which already lines up every label in column A with the correct value one row below
# it (the sheet had been "off by one" between labels and their Portuguese/English
# value pairs from row 13 downward).
$ws.Rows.Item(13).Insert()

# The insert copies column A's formatting into the new A13 cell; the target layout has
# no cell at all in A13, so remove it.
$ws.Range("A13").Clear()

# Row 13 should now hold the professor's name (previously mis-placed two rows below
# "Objetivos:"). Clone formatting from the row below (already correctly styled) so the
# new B13/C13 cells get the right style indexes, then set the value.
$ws.Range("B14").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("B13").Value = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Range("C13").Value = "5840560 - Marco Antonio Carvalho Pereira"

# Fill in the real objectives text under "Objetivos:" (row 10), replacing the
# mis-placed professor name that used to sit there.
$ws.Range("B10").Value = "Fornecer oportunidade de aplicação dos conhecimentos fundamentais da Engenharia de Produção nos diversos sistemas de produção da indústria. Complementação da formação geral curricular. Adaptação psicológica e social do estudante à sua futura atividade profissional."
$ws.Range("C10").Value = "Fornecer oportunidade de aplicação dos conhecimentos fundamentais da Engenharia de Produção nos diversos sistemas de produção da indústria. Complementação da formação geral curricular. Adaptação psicológica e social do estudante à sua futura atividade profissional."

# "Programa resumido:" (row 14, after the insert) gets its real short-syllabus text
# instead of the mis-placed "Semestral".
$ws.Range("B14").Value = "Plano de Trabalho específico. Realização do Estágio. Relatório final e/ou parciais."
$ws.Range("C14").Value = "Plano de Trabalho específico. Realização do Estágio. Relatório final e/ou parciais."

# "Programa:" (row 16, after the insert) gets its real syllabus text instead of the
# mis-placed activation date.
$ws.Range("B16").Value = "Participação do aluno em processo seletivo de empresas ou no setor acadêmico. Estágio realizado sob a supervisão da Escola de Engenharia de Lorena, através do Departamento em Engenharia Química. O conteúdo será estabelecido individualmente no Plano de Trabalho entre o Supervisor do Estágio e o professor orientador, desde que relacionado com as áreas afins da Engenharia de Produção. Apresentação de relatório final e/ou relatórios parciais sobre as atividades desenvolvidas no estágio."
$ws.Range("C16").Value = "Participação do aluno em processo seletivo de empresas ou no setor acadêmico. Estágio realizado sob a supervisão da Escola de Engenharia de Lorena, através do Departamento em Engenharia Química. O conteúdo será estabelecido individualmente no Plano de Trabalho entre o Supervisor do Estágio e o professor orientador, desde que relacionado com as áreas afins da Engenharia de Produção. Apresentação de relatório final e/ou relatórios parciais sobre as atividades desenvolvidas no estágio."

# "Método:" (row 19, after the insert) gets its real method text instead of the
# mis-placed professor name.
$ws.Range("B19").Value = "Supervisão das atividades desenvolvidas pelo aluno durante o estágio."
$ws.Range("C19").Value = "Supervisão das atividades desenvolvidas pelo aluno durante o estágio."

# "Critério:" (row 20, after the insert) gets its real criteria text instead of the
# mis-placed method text.
$ws.Range("B20").Value = "MF = Nota baseada em relatório final e no desempenho no estágio, a ser atribuída pelo professor orientador do estágio."
$ws.Range("C20").Value = "MF = Nota baseada em relatório final e no desempenho no estágio, a ser atribuída pelo professor orientador do estágio."

# "Norma de recuperação:" (row 21, after the insert) gets its real text instead of the
# mis-placed criteria text.
$ws.Range("B21").Value = "Não será oferecida recuperação."
$ws.Range("C21").Value = "Não será oferecida recuperação."

# "Bibliografia:" (row 22, after the insert) gets its real bibliography text instead of
# the mis-placed recovery-policy text.
$ws.Range("B22").Value = "A ser definida com o orientador em função das atividades desenvolvidas no estágio."
$ws.Range("C22").Value = "A ser definida com o orientador em função das atividades desenvolvidas no estágio."
